$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying NATMI TPM recomputation drops all "Target cluster" = ECs rows
# (self/other clusters receiving from ECs-as-target are no longer reported) and
# refreshes every numeric column for the remaining Sending/Target cluster pairs.
# Table shrinks from 9 data rows (A1:T10) to 6 data rows (A1:T7).
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()

$newData = @(
    @{ "A"="ECs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="FAPs"; "E"=3; "F"=1; "G"=35.906979; "H"=107.720937; "I"=0.6107087147789413; "J"=0.6107087147789412; "K"=3; "L"=1; "M"=0.2109236666666666; "N"=0.632771; "O"=0.8951984155054113; "P"=0.8951984155054113; "Q"=7.573631669602999; "R"=68.162685026427; "S"=0.5467054738054544; "T"=0.5467054738054543 },
    @{ "A"="ECs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="MuSCs"; "E"=3; "F"=1; "G"=35.906979; "H"=107.720937; "I"=0.6107087147789413; "J"=0.6107087147789412; "K"=2; "L"=0.6666666666666666; "M"=0.024693; "N"=0.074079; "O"=0.1048015844945887; "P"=0.1048015844945887; "Q"=0.886651032447; "R"=7.979859292023; "S"=0.06400324097348688; "T"=0.06400324097348686 },
    @{ "A"="FAPs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="FAPs"; "E"=3; "F"=1; "G"=17.04862266666667; "H"=51.14586800000001; "I"=0.2899643113254147; "J"=0.2899643113254147; "K"=3; "L"=1; "M"=0.2109236666666666; "N"=0.632771; "O"=0.8951984155054113; "P"=0.8951984155054113; "Q"=3.595958004469778; "R"=32.363622040228; "S"=0.2595755920516291; "T"=0.259575592051629 },
    @{ "A"="FAPs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="MuSCs"; "E"=3; "F"=1; "G"=17.04862266666667; "H"=51.14586800000001; "I"=0.2899643113254147; "J"=0.2899643113254147; "K"=2; "L"=0.6666666666666666; "M"=0.024693; "N"=0.074079; "O"=0.1048015844945887; "P"=0.1048015844945887; "Q"=0.4209816395080002; "R"=3.788834755572001; "S"=0.03038871927378567; "T"=0.03038871927378566 },
    @{ "A"="MuSCs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="FAPs"; "E"=3; "F"=1; "G"=5.839988000000001; "H"=17.519964; "I"=0.09932697389564409; "J"=0.09932697389564407; "K"=3; "L"=1; "M"=0.2109236666666666; "N"=0.632771; "O"=0.8951984155054113; "P"=0.8951984155054113; "Q"=1.231791682249333; "R"=11.086125140244; "S"=0.08891734964832794; "T"=0.08891734964832793 },
    @{ "A"="MuSCs"; "B"="Tgfb1"; "C"="Itgb6"; "D"="MuSCs"; "E"=3; "F"=1; "G"=5.839988000000001; "H"=17.519964; "I"=0.09932697389564409; "J"=0.09932697389564407; "K"=2; "L"=0.6666666666666666; "M"=0.024693; "N"=0.074079; "O"=0.1048015844945887; "P"=0.1048015844945887; "Q"=0.144206823684; "R"=1.297861413156; "S"=0.01040962424731615; "T"=0.01040962424731615 }
)

$r = 2
foreach ($rowData in $newData) {
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$r").Value = $rowData[$col]
    }
    $r = $r + 1
}